$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1169995834814548
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 2.998467759612273

$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 4.327115817150455

$ws.Range("B4").Value = 0.1169995834814548
$ws.Range("C4").Value = 0.04103571897497393
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 1.412515779045154

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.7210945179870265
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("G5").Value = 6.15379541431027

$ws.Range("B6").Value = 0.01253208636536152
$ws.Range("C6").Value = 0.04103571897497393
$ws.Range("D6").Value = 0.1496068669990043
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("G6").Value = 0.7365606309410384

$ws.Range("B7").Value = 0.0006075818656279264
$ws.Range("C7").Value = 0.002658071450198252
$ws.Range("D7").Value = 0.7210945179870265
$ws.Range("E7").Value = 13.86384647080068
$ws.Range("G7").Value = 14.58820664210354

$ws.Range("B8").Value = 3.272327238179451
$ws.Range("C8").Value = 1.626987699542094
$ws.Range("D8").Value = 0.7210945179870265
$ws.Range("E8").Value = 13.86384647080068
$ws.Range("G8").Value = 19.48425592650926

$ws.Range("B9").Value = 1.445647641019636
$ws.Range("C9").Value = 1.626987699542094
$ws.Range("D9").Value = 3.223369029078222
$ws.Range("E9").Value = 0.5333859586016987
$ws.Range("G9").Value = 6.82939032824165

$ws.Range("B10").Value = 1.445647641019636
$ws.Range("C10").Value = 1.626987699542094
$ws.Range("D10").Value = 3.223369029078222
$ws.Range("E10").Value = 13.86384647080068
$ws.Range("G10").Value = 20.15985084044064

$ws.Range("B11").Value = 0.04172184405617529
$ws.Range("C11").Value = 0.002658071450198252
$ws.Range("D11").Value = 189.6080260415259
$ws.Range("E11").Value = 2797.565817734744
$ws.Range("G11").Value = 2987.218223691776

$ws.Range("B12").Value = 3.272327238179451
$ws.Range("C12").Value = 1.626987699542094
$ws.Range("D12").Value = 0.1496068669990043
$ws.Range("E12").Value = 0.5333859586016987
$ws.Range("G12").Value = 5.582307763322248

$ws.Range("B13").Value = 1.445647641019636
$ws.Range("C13").Value = 1.626987699542094
$ws.Range("D13").Value = 3.223369029078222
$ws.Range("E13").Value = 0.5333859586016987
$ws.Range("G13").Value = 6.82939032824165

$ws.Range("B14").Value = 3.272327238179451
$ws.Range("C14").Value = 1.626987699542094
$ws.Range("D14").Value = 0.7210945179870265
$ws.Range("E14").Value = 0.5333859586016987
$ws.Range("G14").Value = 6.15379541431027

$ws.Range("B15").Value = 1.445647641019636
$ws.Range("C15").Value = 1.626987699542094
$ws.Range("D15").Value = 0.7210945179870265
$ws.Range("E15").Value = 13.86384647080068
$ws.Range("G15").Value = 17.65757632934944

$ws.Range("B16").Value = 3.272327238179451
$ws.Range("C16").Value = 1.626987699542094
$ws.Range("D16").Value = 0.1496068669990043
$ws.Range("E16").Value = 0.5333859586016987
$ws.Range("G16").Value = 5.582307763322248

$ws.Range("B17").Value = 3.272327238179451
$ws.Range("C17").Value = 1.626987699542094
$ws.Range("D17").Value = 0.1496068669990043
$ws.Range("E17").Value = 0.5333859586016987
$ws.Range("G17").Value = 5.582307763322248

$ws.Range("B18").Value = 0.01253208636536152
$ws.Range("C18").Value = 0.00006708468553440206
$ws.Range("D18").Value = 0.7210945179870265
$ws.Range("E18").Value = 0.5333859586016987
$ws.Range("G18").Value = 1.267079647639621

$ws.Range("B19").Value = 3.272327238179451
$ws.Range("C19").Value = 1.626987699542094
$ws.Range("D19").Value = 0.1496068669990043
$ws.Range("E19").Value = 0.5333859586016987
$ws.Range("G19").Value = 5.582307763322248

$ws.Range("B20").Value = 1.445647641019636
$ws.Range("C20").Value = 1.626987699542094
$ws.Range("D20").Value = 0.1496068669990043
$ws.Range("E20").Value = 0.5333859586016987
$ws.Range("G20").Value = 3.755628166162433
